$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (13 rows) replacing the original 5 rows
$companies = @(
    @("Simon Property Group Inc", 283),
    @("AMETEK Inc.", 200),
    @("Agilent Technologies Inc", 370),
    @("Autodesk Inc.", 230),
    @("TripAdvisor", 276),
    @("Cabot Oil & Gas", 198),
    @("U.S. Bancorp", 599),
    @("Accenture plc", 455),
    @("Noble Energy Inc", 394),
    @("Lennar Corp.", 465),
    @("PayPal", 169),
    @("Delta Air Lines Inc.", 640),
    @("Pulte Homes Inc.", 631)
)

for ($i = 0; $i -lt $companies.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $companies[$i][0]
    $ws.Cells.Item($row, 2).Value = $companies[$i][1]
}

# Update the chart: style + source data range to cover the new rows
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartStyle = 6

$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,'Sheet1'!`$A`$1:`$A`$13,'Sheet1'!`$B`$1:`$B`$13,1)"
